# connection close fix + random token fix
#
# The author's change added a new "Token new generate when logout" task
# row (row 31, between the existing "Create deployment diagram" row and
# the trailing "Lees opdracht..." summary row 35) to the to-do list table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task row: column C holds the task name, column D holds the "Done" flag (0/1).
$ws.Range("C31").Value = "Token new generate when logout"
$ws.Range("D31").Value = 0

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Range("D32").Select()
